$wb = $excel.ActiveWorkbook

$proyecto = $wb.Worksheets.Item("Proyecto")
$proyecto.Range("E6").Value = 2

$recursos = $wb.Worksheets.Item("Recursos")
$recursos.Activate()
$recursos.Range("B2").Select()

$wb.Worksheets.Item("Totales").Select()
$recursos.Select()
